# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 12357
$ws1.Range("F11").Value = 163
$ws1.Range("F12").Value = 12183
$ws1.Range("F13").Value = 4837
$ws1.Range("F14").Value = 4710
$ws1.Range("F15").Value = 133
$ws1.Range("F16").Value = 65
$ws1.Range("F18").Value = 97
$ws1.Range("F19").Value = 951

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 12357
$ws4.Range("F13").Value = 163
$ws4.Range("F14").Value = 12183
$ws4.Range("F15").Value = 4837
$ws4.Range("F16").Value = 4710
$ws4.Range("F17").Value = 133
$ws4.Range("F18").Value = 65
$ws4.Range("F20").Value = 97
$ws4.Range("F21").Value = 951
